# Fix "excel lock error on Windows os"
# - update the start/end date values
# - move the saved selection off the old (now irrelevant) cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two date cells (serial date numbers, same as stored in the sheet)
$ws.Range("A2").Value = 45748
$ws.Range("B2").Value = 45777

# Make sure the date-range cell style keeps its protection applied
$ws.Range("B2").Locked = $true

# Update the active selection saved with the sheet view
$ws.Activate()
$ws.Range("C4").Select()
